$d = $word.ActiveDocument
$range = $d.Content
$range.Find.Execute("Tanggal : …………………..", $false, $false, $false, $false, $false, $true, 1, $false, "Tanggal : 7 Maret 2025", 2)
